$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.611.37"
$ws.Cells.Item(2, 5).Value = "  -0.03%  "

$ws.Cells.Item(3, 4).Value = "1.596.29"
$ws.Cells.Item(3, 5).Value = "  +0.30%  "

$ws.Cells.Item(4, 5).Value = "  +0.02%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "211.06"
$ws.Cells.Item(5, 5).Value = "  +0.11%  "

$ws.Cells.Item(6, 5).Value = "  +1.39%  "

$ws.Cells.Item(7, 5).Value = "  +0.04%  "

$ws.Cells.Item(8, 5).Value = "  +0.12%  "

$ws.Cells.Item(9, 5).Value = "  -0.86%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.46"
$ws.Cells.Item(10, 5).Value = "  -0.66%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0836"
$ws.Cells.Item(11, 5).Value = "  +0.20%  "

$ws.Cells.Item(12, 4).Value = "1.820.36"
$ws.Cells.Item(12, 5).Value = "  +0.32%  "

$ws.Cells.Item(13, 4).Value = "1.575.34"
$ws.Cells.Item(13, 5).Value = "  -1.14%  "

$ws.Cells.Item(14, 5).Value = "  -0.11%  "

$ws.Cells.Item(15, 5).Value = "  -0.65%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.56"
$ws.Cells.Item(16, 5).Value = "  +0.04%  "

$ws.Cells.Item(17, 4).Value = "26.594.30"
$ws.Cells.Item(17, 5).Value = "  -0.14%  "

$ws.Cells.Item(18, 4).Value = "0.0₃0732"
$ws.Cells.Item(18, 5).Value = "  +0.74%  "

$ws.Cells.Item(19, 5).Value = "  +0.15%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "207.86"
$ws.Cells.Item(20, 5).Value = "  -0.27%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.93"
$ws.Cells.Item(21, 5).Value = "  +3.09%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.24"
$ws.Cells.Item(22, 5).Value = "  +0.08%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.32"
$ws.Cells.Item(23, 5).Value = "  -2.03%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.85"
$ws.Cells.Item(24, 5).Value = "  -0.11%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "145.52"
$ws.Cells.Item(25, 5).Value = "  -0.94%  "

$ws.Cells.Item(26, 5).Value = "  +0.02%  "

$ws.Cells.Item(27, 5).Value = "  -1.61%  "

$ws.Cells.Item(28, 5).Value = "  +0.38%  "

$ws.Cells.Item(29, 5).Value = "  -0.29%  "

$ws.Cells.Item(30, 5).Value = "  -0.90%  "

$ws.Cells.Item(31, 5).Value = "  +0.41%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.22"
$ws.Cells.Item(32, 5).Value = "  -0.02%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.653"
$ws.Cells.Item(33, 5).Value = "  -2.41%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.92"
$ws.Cells.Item(34, 5).Value = "  +0.45%  "

$ws.Cells.Item(35, 4).Value = "1.281.95"
$ws.Cells.Item(35, 5).Value = "  -2.65%  "

$ws.Cells.Item(36, 5).Value = "  +0.79%  "

$ws.Cells.Item(37, 5).Value = "  +0.14%  "

$ws.Cells.Item(38, 5).Value = "  -0.54%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.844"
$ws.Cells.Item(39, 5).Value = "  +1.73%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.44"
$ws.Cells.Item(41, 5).Value = "  +1.40%  "

$ws.Cells.Item(42, 5).Value = "  +1.06%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.785"
$ws.Cells.Item(43, 5).Value = "  -0.64%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "63.94"
$ws.Cells.Item(44, 5).Value = "  +1.60%  "

$ws.Cells.Item(45, 2).Value = "WEMIXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.918"
$ws.Cells.Item(45, 5).Value = "  +10.59%  "

$ws.Cells.Item(46, 2).Value = "RocketPoolETH"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(46, 4).Value = "1.732.62"
$ws.Cells.Item(46, 5).Value = "  +0.34%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "89.55"
$ws.Cells.Item(47, 5).Value = "  -0.36%  "

$ws.Cells.Item(48, 5).Value = "  -1.22%  "

$ws.Cells.Item(49, 5).Value = "  -0.63%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.102"
$ws.Cells.Item(50, 5).Value = "  +3.69%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0506"
$ws.Cells.Item(51, 5).Value = "  -1.33%  "
